$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.966509176771694, 50.06003464001472]"
$ws.Range("T2").Value = "[49.95350288881988, 50.01962549698113]"
$ws.Range("L3").Value = "[49.84682355529793, 49.996310465073705]"
$ws.Range("T3").Value = "[49.95343379039122, 50.054340929615456]"
